# Ajuste de detalles en modulos
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update amounts (Importe) in column B for rows 3, 4 and 6
$ws.Range("B3").Value = 22785.85
$ws.Range("B4").Value = 17720.09
$ws.Range("B6").Value = 14721.55

# Update periodicity (Periodicidad) in column C from "Quincenal" to "Mensual"
# for rows 3 through 6 (this introduces the new "Mensual" shared string)
$ws.Range("C3").Value = "Mensual"
$ws.Range("C4").Value = "Mensual"
$ws.Range("C5").Value = "Mensual"
$ws.Range("C6").Value = "Mensual"

# Reflect the new active cell / selection on the sheet
$ws.Range("E6").Select()
